$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new column before N ---------------
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N; this shifts the old N/O/P/Q columns one
# place to the right (N->O, O->P, P->Q, Q->R) together with their data.
$wsSchedule.Columns("N").Insert() | Out-Null

# Give the freshly inserted column the same width as column M (its left
# neighbour) — this mirrors Excel's own "insert column" width inheritance.
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# --- Make "Repayment Schedule" the active sheet/tab -----------------------
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("T11").Select() | Out-Null
